# ----------------------------------------------------------------------------
# Bosnia Herzegovina Premier Liga -- bases update (25-04-2024 21:26)
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: a handful of fixture rows were re-ordered upstream; swap the
# full data payload (columns B:AB) between each pair while keeping the
# running index in column A attached to its original row. ---

$rowA = $ws.Range("B9:AB9")
$rowB = $ws.Range("B10:AB10")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B29:AB29")
$rowB = $ws.Range("B30:AB30")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B49:AB49")
$rowB = $ws.Range("B50:AB50")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B76:AB76")
$rowB = $ws.Range("B77:AB77")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B87:AB87")
$rowB = $ws.Range("B88:AB88")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B111:AB111")
$rowB = $ws.Range("B112:AB112")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

$rowA = $ws.Range("B122:AB122")
$rowB = $ws.Range("B123:AB123")
$valA = $rowA.Value()
$valB = $rowB.Value()
$rowA.Value = $valB
$rowB.Value = $valA

# --- Step 2: append 6 newly-scraped fixtures (rows 165-170). ---

# Row 165
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(165,1).PasteSpecial(-4122)
$ws.Cells.Item(165,1).Value = 163
$ws.Cells.Item(165,2).Value = 7952755
$ws.Cells.Item(165,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(165,4).Value = 45405.45833333334
$ws.Cells.Item(165,5).Value = "NK Igman Konjic"
$ws.Cells.Item(165,6).Value = "GOSK Gabela"
$ws.Cells.Item(165,7).Value = 0
$ws.Cells.Item(165,8).Value = 0
$ws.Cells.Item(165,9).Value = "D"
$ws.Cells.Item(165,10).Value = 2
$ws.Cells.Item(165,11).Value = 3.3
$ws.Cells.Item(165,12).Value = 3.25
$ws.Cells.Item(165,13).Value = 2
$ws.Cells.Item(165,14).Value = 3.4
$ws.Cells.Item(165,15).Value = 3.2
$ws.Cells.Item(165,16).Value = -0.25
$ws.Cells.Item(165,17).Value = 1.775
$ws.Cells.Item(165,18).Value = 2.025
$ws.Cells.Item(165,19).Value = 2.5
$ws.Cells.Item(165,20).Value = 1.95
$ws.Cells.Item(165,21).Value = 1.85
$ws.Cells.Item(165,22).Value = -1
$ws.Cells.Item(165,23).Value = 2.4
$ws.Cells.Item(165,24).Value = -1
$ws.Cells.Item(165,25).Value = -0.5
$ws.Cells.Item(165,26).Value = 0.5125
$ws.Cells.Item(165,27).Value = -1
$ws.Cells.Item(165,28).Value = 0.8500000000000001
$ws.Cells.Item(165,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 166
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(166,1).PasteSpecial(-4122)
$ws.Cells.Item(166,1).Value = 164
$ws.Cells.Item(166,2).Value = 7952758
$ws.Cells.Item(166,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(166,4).Value = 45405.54166666666
$ws.Cells.Item(166,5).Value = "Velez Mostar"
$ws.Cells.Item(166,6).Value = "Sloga"
$ws.Cells.Item(166,7).Value = 1
$ws.Cells.Item(166,8).Value = 0
$ws.Cells.Item(166,9).Value = "H"
$ws.Cells.Item(166,10).Value = 1.333
$ws.Cells.Item(166,11).Value = 4.5
$ws.Cells.Item(166,12).Value = 7
$ws.Cells.Item(166,13).Value = 1.3
$ws.Cells.Item(166,14).Value = 4.333
$ws.Cells.Item(166,15).Value = 8.5
$ws.Cells.Item(166,16).Value = -1.5
$ws.Cells.Item(166,17).Value = 2
$ws.Cells.Item(166,18).Value = 1.8
$ws.Cells.Item(166,19).Value = 2.5
$ws.Cells.Item(166,20).Value = 1.825
$ws.Cells.Item(166,21).Value = 1.975
$ws.Cells.Item(166,22).Value = 0.3
$ws.Cells.Item(166,23).Value = -1
$ws.Cells.Item(166,24).Value = -1
$ws.Cells.Item(166,25).Value = -1
$ws.Cells.Item(166,26).Value = 0.8
$ws.Cells.Item(166,27).Value = -1
$ws.Cells.Item(166,28).Value = 0.9750000000000001
$ws.Cells.Item(166,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 167
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(167,1).PasteSpecial(-4122)
$ws.Cells.Item(167,1).Value = 165
$ws.Cells.Item(167,2).Value = 7952757
$ws.Cells.Item(167,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(167,4).Value = 45405.64583333334
$ws.Cells.Item(167,5).Value = "Zeljeznicar"
$ws.Cells.Item(167,6).Value = "Siroki Brijeg"
$ws.Cells.Item(167,7).Value = 1
$ws.Cells.Item(167,8).Value = 0
$ws.Cells.Item(167,9).Value = "H"
$ws.Cells.Item(167,10).Value = 1.8
$ws.Cells.Item(167,11).Value = 3.25
$ws.Cells.Item(167,12).Value = 4
$ws.Cells.Item(167,13).Value = 1.4
$ws.Cells.Item(167,14).Value = 4.2
$ws.Cells.Item(167,15).Value = 7
$ws.Cells.Item(167,16).Value = -1.25
$ws.Cells.Item(167,17).Value = 1.95
$ws.Cells.Item(167,18).Value = 1.85
$ws.Cells.Item(167,19).Value = 2.25
$ws.Cells.Item(167,20).Value = 1.85
$ws.Cells.Item(167,21).Value = 1.95
$ws.Cells.Item(167,22).Value = 0.3999999999999999
$ws.Cells.Item(167,23).Value = -1
$ws.Cells.Item(167,24).Value = -1
$ws.Cells.Item(167,25).Value = -0.5
$ws.Cells.Item(167,26).Value = 0.425
$ws.Cells.Item(167,27).Value = -1
$ws.Cells.Item(167,28).Value = 0.95
$ws.Cells.Item(167,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 168
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(168,1).PasteSpecial(-4122)
$ws.Cells.Item(168,1).Value = 166
$ws.Cells.Item(168,2).Value = 7952756
$ws.Cells.Item(168,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(168,4).Value = 45406.45833333334
$ws.Cells.Item(168,5).Value = "Zvijezda 09"
$ws.Cells.Item(168,6).Value = "FK Tuzla City"
$ws.Cells.Item(168,7).Value = 3
$ws.Cells.Item(168,8).Value = 2
$ws.Cells.Item(168,9).Value = "H"
$ws.Cells.Item(168,10).Value = 3.3
$ws.Cells.Item(168,11).Value = 3.4
$ws.Cells.Item(168,12).Value = 1.95
$ws.Cells.Item(168,13).Value = 3
$ws.Cells.Item(168,14).Value = 3.3
$ws.Cells.Item(168,15).Value = 2.15
$ws.Cells.Item(168,16).Value = 0.25
$ws.Cells.Item(168,17).Value = 1.875
$ws.Cells.Item(168,18).Value = 1.925
$ws.Cells.Item(168,19).Value = 2.5
$ws.Cells.Item(168,20).Value = 1.95
$ws.Cells.Item(168,21).Value = 1.85
$ws.Cells.Item(168,22).Value = 2
$ws.Cells.Item(168,23).Value = -1
$ws.Cells.Item(168,24).Value = -1
$ws.Cells.Item(168,25).Value = 0.875
$ws.Cells.Item(168,26).Value = -1
$ws.Cells.Item(168,27).Value = 0.95
$ws.Cells.Item(168,28).Value = -1
$ws.Cells.Item(168,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 169
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(169,1).PasteSpecial(-4122)
$ws.Cells.Item(169,1).Value = 167
$ws.Cells.Item(169,2).Value = 7952759
$ws.Cells.Item(169,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(169,4).Value = 45406.54166666666
$ws.Cells.Item(169,5).Value = "Borac Banja Luka"
$ws.Cells.Item(169,6).Value = "NK Posusje"
$ws.Cells.Item(169,7).Value = 1
$ws.Cells.Item(169,8).Value = 0
$ws.Cells.Item(169,9).Value = "H"
$ws.Cells.Item(169,10).Value = 1.363
$ws.Cells.Item(169,11).Value = 4
$ws.Cells.Item(169,12).Value = 7.5
$ws.Cells.Item(169,13).Value = 1.181
$ws.Cells.Item(169,14).Value = 5.5
$ws.Cells.Item(169,15).Value = 13
$ws.Cells.Item(169,16).Value = -1.75
$ws.Cells.Item(169,17).Value = 1.825
$ws.Cells.Item(169,18).Value = 1.975
$ws.Cells.Item(169,19).Value = 2.75
$ws.Cells.Item(169,20).Value = 2
$ws.Cells.Item(169,21).Value = 1.8
$ws.Cells.Item(169,22).Value = 0.181
$ws.Cells.Item(169,23).Value = -1
$ws.Cells.Item(169,24).Value = -1
$ws.Cells.Item(169,25).Value = -1
$ws.Cells.Item(169,26).Value = 0.9750000000000001
$ws.Cells.Item(169,27).Value = -1
$ws.Cells.Item(169,28).Value = 0.8
$ws.Cells.Item(169,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 170
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(170,1).PasteSpecial(-4122)
$ws.Cells.Item(170,1).Value = 168
$ws.Cells.Item(170,2).Value = 7952460
$ws.Cells.Item(170,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(170,4).Value = 45406.64583333334
$ws.Cells.Item(170,5).Value = "Zrinjski Mostar"
$ws.Cells.Item(170,6).Value = "FK Sarajevo"
$ws.Cells.Item(170,7).Value = 4
$ws.Cells.Item(170,8).Value = 3
$ws.Cells.Item(170,9).Value = "H"
$ws.Cells.Item(170,10).Value = 1.444
$ws.Cells.Item(170,11).Value = 4
$ws.Cells.Item(170,12).Value = 6
$ws.Cells.Item(170,13).Value = 1.444
$ws.Cells.Item(170,14).Value = 3.75
$ws.Cells.Item(170,15).Value = 6.5
$ws.Cells.Item(170,16).Value = -1
$ws.Cells.Item(170,17).Value = 1.775
$ws.Cells.Item(170,18).Value = 2.025
$ws.Cells.Item(170,19).Value = 2.25
$ws.Cells.Item(170,20).Value = 1.925
$ws.Cells.Item(170,21).Value = 1.875
$ws.Cells.Item(170,22).Value = 0.444
$ws.Cells.Item(170,23).Value = -1
$ws.Cells.Item(170,24).Value = -1
$ws.Cells.Item(170,25).Value = 0
$ws.Cells.Item(170,26).Value = 0
$ws.Cells.Item(170,27).Value = 0.925
$ws.Cells.Item(170,28).Value = -1
$ws.Cells.Item(170,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

